$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 10000
$ws.Cells.Item(5, 3).Value = "Tan Nguyen"
$ws.Cells.Item(5, 4).Value = 40465
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat
$ws.Cells.Item(5, 5).Value = "late for meeting"
$ws.Cells.Item(5, 6).Value = "Waiting"

$ws.Range("G4").Select()
